$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# such as "350.40" are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "52.022.72"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "2.869.14"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "350.40"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Value = "112.39"
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "40.24"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "0.0849"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "20.09"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").Value = "3.335.09"
$ws.Range("E15").Value = "  +3.49%  "
$ws.Range("D16").Value = "0.992"
$ws.Range("E16").Value = "  +6.91%  "
$ws.Range("D17").Value = "2.895.61"
$ws.Range("E17").Value = "  +3.54%  "
$ws.Range("D18").Value = "52.070.81"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +4.86%  "
$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "13.58"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").Value = "70.88"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").Value = "270.08"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "2.77"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Value = "26.52"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "10.56"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "38.74"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("D31").Value = "6.25"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").Value = "52.46"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").Value = "5.83"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "0.0454"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0896"
$ws.Range("E35").Value = "  +7.36%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("B37").Value = "Toncoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  -15.90%  "
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  +5.77%  "
$ws.Range("D39").Value = "18.74"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Value = "2.03"
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "122.25"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").Value = "22.46"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").Value = "3.51"
$ws.Range("E46").Value = "  +5.15%  "
$ws.Range("D47").Value = "2.172.86"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("E48").Value = "  +6.42%  "
$ws.Range("D49").Value = "0.238"
$ws.Range("E49").Value = "  +10.10%  "
$ws.Range("D50").Value = "0.960"
$ws.Range("E50").Value = "  +5.81%  "
$ws.Range("D51").Value = "0.0321"
$ws.Range("E51").Value = "  +12.46%  "

# Restore the original (default/general) cell formatting for column D
# now that the text values are safely stored.
$priceRange.ClearFormats()
